$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$values = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 1
    8 = 1
    9 = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 4
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 6
    33 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 1
    38 = 2
    39 = 1
    40 = 0
    41 = 0
    42 = 2
    43 = 1
    44 = 1
    45 = 2
    46 = 1
    47 = 3
    48 = 2
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 2
    54 = 0
    55 = 0
    56 = 1
    57 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $values[$row]
}
